$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.382.71'
$ws.Range("E2").Value = '  -2.26%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.844.16'
$ws.Range("E3").Value = '  -1.85%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '260.41'
$ws.Range("E5").Value = '  -7.45%  '

$ws.Range("E6").Value = '  +0.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5249'
$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3236'
$ws.Range("E8").Value = '  -8.28%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06751'
$ws.Range("E9").Value = '  -4.55%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.91'
$ws.Range("E10").Value = '  -6.90%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7725'
$ws.Range("E11").Value = '  -5.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07689'
$ws.Range("E12").Value = '  -1.47%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.814.52'
$ws.Range("E13").Value = '  -3.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '89.02'
$ws.Range("E14").Value = '  -1.53%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.033'
$ws.Range("E15").Value = '  -3.72%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.03%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.15'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007885'
$ws.Range("E19").Value = '  -3.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.415.77'
$ws.Range("E20").Value = '  -2.24%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.080.43'
$ws.Range("E21").Value = '  -2.23%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.537'
$ws.Range("E22").Value = '  -5.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.473'
$ws.Range("E23").Value = '  -7.19%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.929'
$ws.Range("E24").Value = '  -5.20%  '

$ws.Range("E25").Value = '  -2.78%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '144.27'
$ws.Range("E26").Value = '  -1.75%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.647'
$ws.Range("E27").Value = '  -1.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.94'
$ws.Range("E28").Value = '  -3.75%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.53'
$ws.Range("E29").Value = '  -1.94%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.195'
$ws.Range("E30").Value = '  -4.56%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08802'
$ws.Range("E31").Value = '  -0.91%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.138'
$ws.Range("E32").Value = '  -5.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04846'
$ws.Range("E33").Value = '  -1.10%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.134'
$ws.Range("E34").Value = '  -3.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.853'
$ws.Range("E35").Value = '  -1.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6825'
$ws.Range("E36").Value = '  -8.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.112'
$ws.Range("E37").Value = '  -5.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01794'
$ws.Range("E38").Value = '  -4.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.214'
$ws.Range("E39").Value = '  -8.23%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4922'
$ws.Range("E40").Value = '  -7.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '112.77'
$ws.Range("E41").Value = '  -3.70%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8994'
$ws.Range("E42").Value = '  -8.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.174'
$ws.Range("E43").Value = '  -2.32%  '

$ws.Range("E44").Value = '  +0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.760'
$ws.Range("E45").Value = '  -5.00%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4195'
$ws.Range("E46").Value = '  -9.10%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1263'
$ws.Range("E47").Value = '  -7.77%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.099'
$ws.Range("E48").Value = '  -4.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05882'
$ws.Range("E49").Value = '  -1.09%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '35.40'
$ws.Range("E50").Value = '  -3.82%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.35'
$ws.Range("E51").Value = '  -4.17%  '
